$d = $word.ActiveDocument

# Locate the paragraph that ends with "...is not written out by Qualtrics.  "
# (the AllOneSide paragraph) so we can insert the new timeDelay paragraph
# immediately after it, before the "Amazon Mechanical Turk" heading.
$anchorIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*is not written out by Qualtrics*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not locate the AllOneSide / 'is not written out by Qualtrics' paragraph."
}

$anchorPara = $d.Paragraphs.Item($anchorIndex)

# Insert a brand-new (empty) paragraph right after the AllOneSide paragraph.
$anchorPara.Range.InsertParagraphAfter()

# The new empty paragraph is now at index (anchorIndex + 1); fill it in with
# the new "timeDelay" paragraph text.
$newPara = $d.Paragraphs.Item($anchorIndex + 1)
$newPara.Range.Text = "A timing variable, timeDelay, is used to set a slight delay before each image appears.  For images with high levels of similarity, this can make it easier for subjects to realize that there is a new set of images.  The default is 250ms.  "
